$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and one row swap)
# Force text number-format before assigning, then restore default style,
# so values such as "509.24" or "314.00" are kept as exact text and not
# coerced into floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.446.36'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.421.08'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '509.24'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.02'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.67%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.445.63'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0971'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.62'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.856.43'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.18%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.366.72'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.91'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.71%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.399.38'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.29'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '314.00'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.06'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.96%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.992'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.382'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.56'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.12%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.55'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0730'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.68'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.01'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.12'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.02'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.91%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.85'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.57%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.72'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.804'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '132.12'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.68%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('B44').Style = 'Normal'

$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C44').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.01'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.19%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('B45').Style = 'Normal'

$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C45').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.38'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '254.06'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.569'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0914'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0491'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.21'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.77%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.71%  '
$ws.Range('E51').Style = 'Normal'

